$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 258.22223
$ws.Range("I11").Value = 258.22223
$ws.Range("K11").Value = 258.22223
$ws.Range("M11").Value = -118.22223
$ws.Range("H28").Value = 12870.471
$ws.Range("J28").Value = 21199.375
$ws.Range("L28").Value = 21199.375
$ws.Range("N28").Value = -22169.375
$ws.Range("H94").Value = 10153.111
$ws.Range("I94").Value = 10153.111
$ws.Range("K94").Value = 10153.111
$ws.Range("M94").Value = -9702.111000000001
$ws.Range("H101").Value = 12692.75
$ws.Range("I101").Value = 11591.8
$ws.Range("J101").Value = 14527.667
$ws.Range("K101").Value = 34775.39999999999
$ws.Range("L101").Value = 43583.001
$ws.Range("M101").Value = -33153.39999999999
$ws.Range("N101").Value = -46827.001
$ws.Range("H113").Value = 4895.7915
$ws.Range("J113").Value = 3518.6667
$ws.Range("L113").Value = 3518.6667
$ws.Range("N113").Value = -10026.6667
$ws.Range("H125").Value = 4950.2856
$ws.Range("I125").Value = 4144
$ws.Range("K125").Value = 37296
$ws.Range("M125").Value = -34836
$ws.Range("H127").Value = 1788.6
$ws.Range("I127").Value = 1981
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 5943
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = -983
$ws.Range("N127").Value = -14420
$ws.Range("H138").Value = 3138.56
$ws.Range("J138").Value = 3416
$ws.Range("L138").Value = 10248
$ws.Range("N138").Value = -20528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3492.8333
$ws.Range("I61").Value = 3199.3333
$ws.Range("K61").Value = 3199.3333
$ws.Range("M61").Value = -2987.3333
$ws.Range("H74").Value = 2072.7805
$ws.Range("I74").Value = 1710.9166
$ws.Range("K74").Value = 1710.9166
$ws.Range("M74").Value = -836.9166
$ws.Range("H77").Value = 2072.7805
$ws.Range("I77").Value = 1710.9166
$ws.Range("K77").Value = 8554.583000000001
$ws.Range("M77").Value = -4186.583000000001
$ws.Range("H102").Value = 1799.9
$ws.Range("I102").Value = 1821.7059
$ws.Range("K102").Value = 1821.7059
$ws.Range("M102").Value = -199.7058999999999
$ws.Range("H110").Value = 2272.2
$ws.Range("I110").Value = 1961.7142
$ws.Range("J110").Value = 2996.6667
$ws.Range("K110").Value = 1961.7142
$ws.Range("L110").Value = 2996.6667
$ws.Range("M110").Value = 83.28580000000011
$ws.Range("N110").Value = -7086.6667
$ws.Range("H132").Value = 1932.2449
$ws.Range("I132").Value = 1186
$ws.Range("K132").Value = 3558
$ws.Range("M132").Value = -1028
$ws.Range("H135").Value = 67499.5
$ws.Range("J135").Value = 67499.5
$ws.Range("L135").Value = 67499.5
$ws.Range("N135").Value = -77639.5
$ws.Range("H136").Value = 3492.8333
$ws.Range("I136").Value = 3199.3333
$ws.Range("K136").Value = 9597.999899999999
$ws.Range("M136").Value = -7047.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 58750
$ws.Range("J132").Value = 58750
$ws.Range("L132").Value = 58750
$ws.Range("N132").Value = -68870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2500
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3574
$ws.Range("H103").Value = 7371.5
$ws.Range("I103").Value = 7371.5
$ws.Range("K103").Value = 7371.5
$ws.Range("M103").Value = -6199.5
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18790074
$ws.Range("I4").Value = 1979348
$ws.Range("K4").Value = 5938044
$ws.Range("M4").Value = -5937932
$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -10068
$ws.Range("H80").Value = 12703.947
$ws.Range("J80").Value = 7636.4
$ws.Range("L80").Value = 22909.2
$ws.Range("N80").Value = -24781.2
$ws.Range("H83").Value = 12703.947
$ws.Range("J83").Value = 7636.4
$ws.Range("L83").Value = 68727.59999999999
$ws.Range("N83").Value = -78087.59999999999
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -2400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 69995
$ws.Range("J109").Value = 69995
$ws.Range("L109").Value = 69995
$ws.Range("N109").Value = -72075
$ws.Range("H132").Value = 2916.8708
$ws.Range("I132").Value = 2949.75
$ws.Range("J132").Value = 2804.1428
$ws.Range("K132").Value = 8849.25
$ws.Range("L132").Value = 8412.428400000001
$ws.Range("M132").Value = -6319.25
$ws.Range("N132").Value = -13472.4284
$ws.Range("H137").Value = 66000
$ws.Range("J137").Value = 66000
$ws.Range("L137").Value = 66000
$ws.Range("N137").Value = -76200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3799.5625
$ws.Range("I7").Value = 4025.0833
$ws.Range("J7").Value = 3123
$ws.Range("K7").Value = 4025.0833
$ws.Range("L7").Value = 3123
$ws.Range("M7").Value = -3913.0833
$ws.Range("N7").Value = -3347
$ws.Range("H43").Value = 5696583.5
$ws.Range("J43").Value = 5696583.5
$ws.Range("L43").Value = 5696583.5
$ws.Range("N43").Value = -5696969.5
$ws.Range("H61").Value = 9980.546
$ws.Range("I61").Value = 9977.556
$ws.Range("K61").Value = 9977.556
$ws.Range("M61").Value = -9775.556
$ws.Range("H113").Value = 9980.546
$ws.Range("I113").Value = 9977.556
$ws.Range("K113").Value = 9977.556
$ws.Range("M113").Value = -7807.556
$ws.Range("H121").Value = 110000.5
$ws.Range("J121").Value = 110000.5
$ws.Range("L121").Value = 110000.5
$ws.Range("N121").Value = -113494.5
$ws.Range("H126").Value = 3799.5625
$ws.Range("I126").Value = 4025.0833
$ws.Range("J126").Value = 3123
$ws.Range("K126").Value = 12075.2499
$ws.Range("L126").Value = 9369
$ws.Range("M126").Value = -9605.249899999999
$ws.Range("N126").Value = -14309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3333.3333
$ws.Range("I81").Value = 3500
$ws.Range("K81").Value = 7000
$ws.Range("M81").Value = -5939
$ws.Range("H84").Value = 3333.3333
$ws.Range("I84").Value = 3500
$ws.Range("K84").Value = 35000
$ws.Range("M84").Value = -29696
$ws.Range("H113").Value = 1900.2222
$ws.Range("I113").Value = 1874.5
$ws.Range("J113").Value = 1920.8
$ws.Range("K113").Value = 5623.5
$ws.Range("L113").Value = 5762.4
$ws.Range("M113").Value = -3453.5
$ws.Range("N113").Value = -10102.4
$ws.Range("H122").Value = 3950.4119
$ws.Range("I122").Value = 3649.5833
$ws.Range("J122").Value = 4672.4
$ws.Range("K122").Value = 10948.7499
$ws.Range("L122").Value = 14017.2
$ws.Range("M122").Value = -8498.749899999999
$ws.Range("N122").Value = -18917.2
$ws.Range("H132").Value = 4017.8462
$ws.Range("I132").Value = 3848.6365
$ws.Range("K132").Value = 11545.9095
$ws.Range("M132").Value = -9015.9095
$ws.Range("H140").Value = 74999
$ws.Range("J140").Value = 74999
$ws.Range("L140").Value = 74999
$ws.Range("N140").Value = -85359
$ws.Range("H141").Value = 106000
$ws.Range("J141").Value = 106000
$ws.Range("L141").Value = 106000
$ws.Range("N141").Value = -116360
